$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 85322.09149999999
$ws.Range("C2").Value = 5539.9085
$ws.Range("D2").Value = 65000
$ws.Range("E2").Value = 14782.18299999999

$ws.Range("B3").Value = 66643.49400000001
$ws.Range("C3").Value = 5343.505999999999
$ws.Range("D3").Value = 55000
$ws.Range("E3").Value = 6299.988000000005

$ws.Range("B4").Value = 63850.3595
$ws.Range("C4").Value = 5250.640500000001
$ws.Range("D4").Value = 45000
$ws.Range("E4").Value = 13599.719

$ws.Range("B5").Value = 61979.8705
$ws.Range("C5").Value = 5194.129499999999
$ws.Range("D5").Value = 45000
$ws.Range("E5").Value = 11785.74099999999

$ws.Range("B6").Value = 62565.603
$ws.Range("C6").Value = 5143.397
$ws.Range("E6").Value = 24922.20600000001

$ws.Range("B7").Value = 74491.1695
$ws.Range("C7").Value = 5188.8305
$ws.Range("E7").Value = 46802.33900000001

$ws.Range("B8").Value = 60588.8635
$ws.Range("C8").Value = 5628.136500000001
$ws.Range("E8").Value = 32460.727

$ws.Range("B9").Value = 67205.92599999999
$ws.Range("C9").Value = 6549.074000000001
$ws.Range("E9").Value = 38156.85199999999

$ws.Range("B10").Value = 75170.359
$ws.Range("C10").Value = 7809.641
$ws.Range("D10").Value = 22500
$ws.Range("E10").Value = 44860.71799999999

$ws.Range("B11").Value = 82295.1295
$ws.Range("C11").Value = 12977.8705
$ws.Range("D11").Value = 32500
$ws.Range("E11").Value = 36817.25899999999

$ws.Range("B12").Value = 83060.50599999999
$ws.Range("C12").Value = 14879.494
$ws.Range("D12").Value = 55000
$ws.Range("E12").Value = 13181.01199999999

$ws.Range("B13").Value = 79680.716
$ws.Range("C13").Value = 15334.284
$ws.Range("E13").Value = -653.5679999999993

$ws.Range("B14").Value = 87679.5435
$ws.Range("C14").Value = 15538.4565
$ws.Range("E14").Value = 7141.087

$ws.Range("B15").Value = 85245.04149999999
$ws.Range("C15").Value = 15473.9585
$ws.Range("E15").Value = 4771.082999999984

$ws.Range("B16").Value = 84102.6715
$ws.Range("C16").Value = 15459.3285
$ws.Range("E16").Value = 3643.342999999993

$ws.Range("B17").Value = 82366.4325
$ws.Range("C17").Value = 15733.5675
$ws.Range("E17").Value = 1632.864999999991

$ws.Range("B18").Value = 78912.772
$ws.Range("C18").Value = 15999.228
$ws.Range("E18").Value = -2086.456000000006

$ws.Range("B19").Value = 77523.92999999999
$ws.Range("C19").Value = 15687.07
$ws.Range("E19").Value = -3163.140000000007

$ws.Range("B20").Value = 79823.067
$ws.Range("C20").Value = 15318.933
$ws.Range("E20").Value = -495.8660000000091

$ws.Range("B21").Value = 77799.3585
$ws.Range("C21").Value = 13434.6415
$ws.Range("E21").Value = -635.2829999999958

$ws.Range("B22").Value = 77147.185
$ws.Range("C22").Value = 11609.815
$ws.Range("E22").Value = 537.3699999999953

$ws.Range("B23").Value = 75950.549
$ws.Range("C23").Value = 9103.451000000001
$ws.Range("E23").Value = 1847.097999999998

$ws.Range("B24").Value = 73885.5585
$ws.Range("C24").Value = 6290.4415
$ws.Range("E24").Value = 2595.116999999998

$ws.Range("B25").Value = 71567.1335
$ws.Range("C25").Value = 5217.8665
$ws.Range("E25").Value = 1349.266999999993
